$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-4: Sending cluster = ECs; Row 5-7: FAPs; Row 8-10: sCs
# Target cluster cycles ECs, FAPs, sCs within each block

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Inhbb"
$ws.Cells.Item(2, 3).Value = "Acvr1b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 3.345805333333333
$ws.Cells.Item(2, 8).Value = 10.037416
$ws.Cells.Item(2, 9).Value = 0.3489465220682754
$ws.Cells.Item(2, 10).Value = 0.3489465220682754
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.021200333333334
$ws.Cells.Item(2, 14).Value = 12.063601
$ws.Cells.Item(2, 15).Value = 0.389801966361343
$ws.Cells.Item(2, 16).Value = 0.389801966361343
$ws.Cells.Item(2, 17).Value = 13.45415352166844
$ws.Cells.Item(2, 18).Value = 121.087381695016
$ws.Cells.Item(2, 19).Value = 0.1360200404571655
$ws.Cells.Item(2, 20).Value = 0.1360200404571655

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Inhbb"
$ws.Cells.Item(3, 3).Value = "Acvr1b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 3.345805333333333
$ws.Cells.Item(3, 8).Value = 10.037416
$ws.Cells.Item(3, 9).Value = 0.3489465220682754
$ws.Cells.Item(3, 10).Value = 0.3489465220682754
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.562995333333333
$ws.Cells.Item(3, 14).Value = 10.688986
$ws.Cells.Item(3, 15).Value = 0.3453850770768087
$ws.Cells.Item(3, 16).Value = 0.3453850770768087
$ws.Cells.Item(3, 17).Value = 11.92108878890844
$ws.Cells.Item(3, 18).Value = 107.289799100176
$ws.Cells.Item(3, 19).Value = 0.1205209214202356
$ws.Cells.Item(3, 20).Value = 0.1205209214202356

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Inhbb"
$ws.Cells.Item(4, 3).Value = "Acvr1b"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 3.345805333333333
$ws.Cells.Item(4, 8).Value = 10.037416
$ws.Cells.Item(4, 9).Value = 0.3489465220682754
$ws.Cells.Item(4, 10).Value = 0.3489465220682754
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.731812666666666
$ws.Cells.Item(4, 14).Value = 8.195438
$ws.Cells.Item(4, 15).Value = 0.2648129565618484
$ws.Cells.Item(4, 16).Value = 0.2648129565618484
$ws.Cells.Item(4, 17).Value = 9.140113389800886
$ws.Cells.Item(4, 18).Value = 82.26102050820798
$ws.Cells.Item(4, 19).Value = 0.0924055601908743
$ws.Cells.Item(4, 20).Value = 0.0924055601908743

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Inhbb"
$ws.Cells.Item(5, 3).Value = "Acvr1b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.480061666666667
$ws.Cells.Item(5, 8).Value = 16.440185
$ws.Cells.Item(5, 9).Value = 0.5715360783999618
$ws.Cells.Item(5, 10).Value = 0.5715360783999618
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.021200333333334
$ws.Cells.Item(5, 14).Value = 12.063601
$ws.Cells.Item(5, 15).Value = 0.389801966361343
$ws.Cells.Item(5, 16).Value = 0.389801966361343
$ws.Cells.Item(5, 17).Value = 22.03642580068723
$ws.Cells.Item(5, 18).Value = 198.327832206185
$ws.Cells.Item(5, 19).Value = 0.2227858872067558
$ws.Cells.Item(5, 20).Value = 0.2227858872067558

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Inhbb"
$ws.Cells.Item(6, 3).Value = "Acvr1b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.480061666666667
$ws.Cells.Item(6, 8).Value = 16.440185
$ws.Cells.Item(6, 9).Value = 0.5715360783999618
$ws.Cells.Item(6, 10).Value = 0.5715360783999618
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.562995333333333
$ws.Cells.Item(6, 14).Value = 10.688986
$ws.Cells.Item(6, 15).Value = 0.3453850770768087
$ws.Cells.Item(6, 16).Value = 0.3453850770768087
$ws.Cells.Item(6, 17).Value = 19.52543414471222
$ws.Cells.Item(6, 18).Value = 175.72890730241
$ws.Cells.Item(6, 19).Value = 0.1974000324903478
$ws.Cells.Item(6, 20).Value = 0.1974000324903478

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Inhbb"
$ws.Cells.Item(7, 3).Value = "Acvr1b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.480061666666667
$ws.Cells.Item(7, 8).Value = 16.440185
$ws.Cells.Item(7, 9).Value = 0.5715360783999618
$ws.Cells.Item(7, 10).Value = 0.5715360783999618
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.731812666666666
$ws.Cells.Item(7, 14).Value = 8.195438
$ws.Cells.Item(7, 15).Value = 0.2648129565618484
$ws.Cells.Item(7, 16).Value = 0.2648129565618484
$ws.Cells.Item(7, 17).Value = 14.97050187511444
$ws.Cells.Item(7, 18).Value = 134.73451687603
$ws.Cells.Item(7, 19).Value = 0.1513501587028583
$ws.Cells.Item(7, 20).Value = 0.1513501587028583

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Inhbb"
$ws.Cells.Item(8, 3).Value = "Acvr1b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.7624369999999999
$ws.Cells.Item(8, 8).Value = 2.287311
$ws.Cells.Item(8, 9).Value = 0.07951739953176286
$ws.Cells.Item(8, 10).Value = 0.07951739953176286
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.021200333333334
$ws.Cells.Item(8, 14).Value = 12.063601
$ws.Cells.Item(8, 15).Value = 0.389801966361343
$ws.Cells.Item(8, 16).Value = 0.389801966361343
$ws.Cells.Item(8, 17).Value = 3.065911918545666
$ws.Cells.Item(8, 18).Value = 27.593207266911
$ws.Cells.Item(8, 19).Value = 0.0309960386974217
$ws.Cells.Item(8, 20).Value = 0.03099603869742169

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Inhbb"
$ws.Cells.Item(9, 3).Value = "Acvr1b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.7624369999999999
$ws.Cells.Item(9, 8).Value = 2.287311
$ws.Cells.Item(9, 9).Value = 0.07951739953176286
$ws.Cells.Item(9, 10).Value = 0.07951739953176286
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.562995333333333
$ws.Cells.Item(9, 14).Value = 10.688986
$ws.Cells.Item(9, 15).Value = 0.3453850770768087
$ws.Cells.Item(9, 16).Value = 0.3453850770768087
$ws.Cells.Item(9, 17).Value = 2.716559472960666
$ws.Cells.Item(9, 18).Value = 24.449035256646
$ws.Cells.Item(9, 19).Value = 0.0274641231662253
$ws.Cells.Item(9, 20).Value = 0.0274641231662253

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Inhbb"
$ws.Cells.Item(10, 3).Value = "Acvr1b"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.7624369999999999
$ws.Cells.Item(10, 8).Value = 2.287311
$ws.Cells.Item(10, 9).Value = 0.07951739953176286
$ws.Cells.Item(10, 10).Value = 0.07951739953176286
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.731812666666666
$ws.Cells.Item(10, 14).Value = 8.195438
$ws.Cells.Item(10, 15).Value = 0.2648129565618484
$ws.Cells.Item(10, 16).Value = 0.2648129565618484
$ws.Cells.Item(10, 17).Value = 2.082835054135333
$ws.Cells.Item(10, 18).Value = 18.745515487218
$ws.Cells.Item(10, 19).Value = 0.02105723766811586
$ws.Cells.Item(10, 20).Value = 0.02105723766811586
